$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFCB5")

# Make IFCB5 the active sheet (this flips tabSelected from IFCB101 -> IFCB5
# and updates bookViews/workbookView activeTab accordingly).
$ws.Activate() | Out-Null

# New bead-voltage entry row (row 3).
$ws.Range("A3").Value = "IFCB5_2015_285_180748"

# B3 / C3 need the same date / time number formats as B2 / C2 (styles 8 & 9).
# Copy formats first so we land on the existing shared styles instead of
# Excel minting brand-new numFmt/style entries, then set the values.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = 42289

$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 0.42152777777777778

$ws.Range("E3").Value = 1.6
$ws.Range("F3").Value = 9.8000000000000007
$ws.Range("G3").Value = 7.7
$ws.Range("I3").Value = 1.1000000000000001
$ws.Range("J3").Value = 24.8
$ws.Range("K3").Value = 41.5
$ws.Range("M3").Value = "9um beads in SW after changing delay to about 270us guess that's not what it was before. Could be new board?"

# Selection lands on A4, the row right after the new entry.
# (IFCB101, the previously-active sheet, keeps its old E5 selection
# untouched -- it only loses tabSelected as a side effect of IFCB5
# becoming the active sheet above.)
$ws.Range("A4").Select() | Out-Null
